$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1551020408163265
$ws.Range("C2").Value = 0.6408163265306123
$ws.Range("J2").Value = 0.01224489795918367
$ws.Range("P2").Value = 0.09795918367346938
$ws.Range("S2").Value = 0.09387755102040816
$ws.Range("C3").Value = 0.01851851851851852
$ws.Range("J3").Value = 0.03703703703703703
$ws.Range("P3").Value = 0.7530864197530864
$ws.Range("S3").Value = 0.191358024691358
$ws.Range("P4").Value = 0.7105263157894737
$ws.Range("S4").Value = 0.2894736842105263
$ws.Range("B6").Value = 0.04587155963302753
$ws.Range("D6").Value = 0.009174311926605505
$ws.Range("F6").Value = 0.03669724770642202
$ws.Range("J6").Value = 0.2339449541284404
$ws.Range("O6").Value = 0.009174311926605505
$ws.Range("Q6").Value = 0.1972477064220184
$ws.Range("R6").Value = 0.1055045871559633
$ws.Range("S6").Value = 0.3623853211009174
$ws.Range("B7").Value = 0.06993006993006994
$ws.Range("D7").Value = 0.006993006993006993
$ws.Range("F7").Value = 0.06993006993006994
$ws.Range("J7").Value = 0.1258741258741259
$ws.Range("O7").Value = 0.02797202797202797
$ws.Range("Q7").Value = 0.1818181818181818
$ws.Range("R7").Value = 0.1188811188811189
$ws.Range("S7").Value = 0.3986013986013986
$ws.Range("B8").Value = 0.08064516129032258
$ws.Range("D8").Value = 0.01382488479262673
$ws.Range("F8").Value = 0.05529953917050692
$ws.Range("J8").Value = 0.1336405529953917
$ws.Range("O8").Value = 0.01152073732718894
$ws.Range("Q8").Value = 0.1866359447004608
$ws.Range("R8").Value = 0.1129032258064516
$ws.Range("S8").Value = 0.4055299539170507
$ws.Range("B9").Value = 0.1019607843137255
$ws.Range("D9").Value = 0.007843137254901961
$ws.Range("F9").Value = 0.05490196078431372
$ws.Range("J9").Value = 0.1411764705882353
$ws.Range("O9").Value = 0.01176470588235294
$ws.Range("Q9").Value = 0.1529411764705882
$ws.Range("R9").Value = 0.1176470588235294
$ws.Range("S9").Value = 0.4117647058823529
$ws.Range("B10").Value = 0.08553546592489569
$ws.Range("D10").Value = 0.01947148817802503
$ws.Range("F10").Value = 0.06884561891515995
$ws.Range("J10").Value = 0.1363004172461753
$ws.Range("O10").Value = 0.01390820584144645
$ws.Range("Q10").Value = 0.217663421418637
$ws.Range("R10").Value = 0.10778859527121
$ws.Range("S10").Value = 0.3504867872044506
$ws.Range("G11").Value = 0.16289592760181
$ws.Range("J11").Value = 0.06787330316742081
$ws.Range("K11").Value = 0.2126696832579185
$ws.Range("L11").Value = 0.5475113122171946
$ws.Range("S11").Value = 0.009049773755656109
$ws.Range("G12").Value = 0.6535433070866141
$ws.Range("J12").Value = 0.2440944881889764
$ws.Range("K12").Value = 0.01574803149606299
$ws.Range("L12").Value = 0.05511811023622047
$ws.Range("S12").Value = 0.03149606299212598
$ws.Range("G13").Value = 0.7105263157894737
$ws.Range("J13").Value = 0.2368421052631579
$ws.Range("S13").Value = 0.05263157894736842
$ws.Range("F15").Value = 0.01941747572815534
$ws.Range("H15").Value = 0.1747572815533981
$ws.Range("I15").Value = 0.07766990291262135
$ws.Range("J15").Value = 0.3300970873786408
$ws.Range("K15").Value = 0.06310679611650485
$ws.Range("M15").Value = 0.01941747572815534
$ws.Range("N15").Value = 0.004854368932038835
$ws.Range("O15").Value = 0.08737864077669903
$ws.Range("S15").Value = 0.2233009708737864
$ws.Range("F16").Value = 0.01219512195121951
$ws.Range("H16").Value = 0.1463414634146341
$ws.Range("I16").Value = 0.1158536585365854
$ws.Range("J16").Value = 0.4878048780487805
$ws.Range("K16").Value = 0.07926829268292683
$ws.Range("M16").Value = 0.01829268292682927
$ws.Range("N16").Value = 0.006097560975609756
$ws.Range("O16").Value = 0.06097560975609756
$ws.Range("S16").Value = 0.07317073170731707
$ws.Range("F17").Value = 0.01803607214428858
$ws.Range("H17").Value = 0.1983967935871744
$ws.Range("I17").Value = 0.08617234468937876
$ws.Range("J17").Value = 0.468937875751503
$ws.Range("K17").Value = 0.0501002004008016
$ws.Range("M17").Value = 0.01402805611222445
$ws.Range("O17").Value = 0.0561122244488978
$ws.Range("S17").Value = 0.1082164328657315
$ws.Range("F18").Value = 0.01470588235294118
$ws.Range("H18").Value = 0.1397058823529412
$ws.Range("I18").Value = 0.09191176470588236
$ws.Range("J18").Value = 0.5110294117647058
$ws.Range("K18").Value = 0.04779411764705882
$ws.Range("M18").Value = 0.02205882352941177
$ws.Range("O18").Value = 0.05514705882352941
$ws.Range("S18").Value = 0.1176470588235294
$ws.Range("F19").Value = 0.01250977326035966
$ws.Range("H19").Value = 0.1837372947615324
$ws.Range("I19").Value = 0.1196247068021892
$ws.Range("J19").Value = 0.4073494917904613
$ws.Range("K19").Value = 0.07896794370602032
$ws.Range("M19").Value = 0.01407349491790461
$ws.Range("O19").Value = 0.06020328381548085
$ws.Range("S19").Value = 0.1235340109460516
